$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.180087685585022
$ws.Range("B1").Value = 2.417452573776245
$ws.Range("D1").Value = 2.331658601760864
$ws.Range("E1").Value = 1.193108081817627
